# Auto-generated edit script: updates cryptocurrency Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.502.55'
$ws.Range("E2").Value = '  +0.87%  '
$ws.Range("D3").Value = '2.603.08'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").Value = '''539.19'
$ws.Range("E5").Value = '  +3.43%  '
$ws.Range("D6").Value = '''141.58'
$ws.Range("E6").Value = '  +1.62%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '''0.565'
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '''6.51'
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("D10").Value = '''0.103'
$ws.Range("E10").Value = '  +1.61%  '
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("E12").Value = '  -0.93%  '
$ws.Range("D13").Value = '3.060.90'
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").Value = '59.404.25'
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("D15").Value = '''20.82'
$ws.Range("E15").Value = '  +1.16%  '
$ws.Range("D16").Value = '2.626.58'
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").Value = '''341.38'
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").Value = '''4.38'
$ws.Range("E19").Value = '  +1.56%  '
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("D21").Value = '''6.33'
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '''67.29'
$ws.Range("E23").Value = '  +1.73%  '
$ws.Range("D24").Value = '''0.409'
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("E25").Value = '  -1.45%  '
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").Value = '''7.20'
$ws.Range("E27").Value = '  +2.53%  '
$ws.Range("D28").Value = '0.0₃0746'
$ws.Range("E28").Value = '  +2.92%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  +6.09%  '
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("D32").Value = '''18.83'
$ws.Range("E32").Value = '  +0.71%  '
$ws.Range("D33").Value = '''149.90'
$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("D34").Value = '''3.98'
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("D36").Value = '''0.846'
$ws.Range("E36").Value = '  +3.50%  '
$ws.Range("E37").Value = '  -0.68%  '
$ws.Range("D38").Value = '''0.827'
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("D39").Value = '''3.54'
$ws.Range("E39").Value = '  +0.41%  '
$ws.Range("D40").Value = '''0.999'
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("D41").Value = '''273.23'
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("E42").Value = '  +1.40%  '
$ws.Range("D43").Value = '''10.74'
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("D44").Value = '''0.0951'
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("E45").Value = '  +1.40%  '
$ws.Range("D46").Value = '''18.56'
$ws.Range("E46").Value = '  +3.85%  '
$ws.Range("E47").Value = '  +1.49%  '
$ws.Range("D48").Value = '1.938.72'
$ws.Range("E48").Value = '  -1.56%  '
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("D50").Value = '''111.92'
$ws.Range("E50").Value = '  -1.50%  '
